$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 112643.03
$ws.Range("I15").Value = 112643.03
$ws.Range("K15").Value = 337929.09
$ws.Range("M15").Value = -337760.09
$ws.Range("H53").Value = 71694.64
$ws.Range("I53").Value = 142935.58
$ws.Range("J53").Value = 453.7143
$ws.Range("K53").Value = 142935.58
$ws.Range("L53").Value = 453.7143
$ws.Range("M53").Value = -142298.58
$ws.Range("N53").Value = -1727.7143
$ws.Range("H106").Value = 2780.2727
$ws.Range("I106").Value = 2287.111
$ws.Range("K106").Value = 2287.111
$ws.Range("M106").Value = -1656.111
$ws.Range("H107").Value = 1558.5778
$ws.Range("I107").Value = 1725.7028
$ws.Range("J107").Value = 785.625
$ws.Range("K107").Value = 1725.7028
$ws.Range("L107").Value = 785.625
$ws.Range("M107").Value = 194.2972
$ws.Range("N107").Value = -4625.625
$ws.Range("H132").Value = 260524.95
$ws.Range("I132").Value = 303924.16
$ws.Range("J132").Value = 5554.625
$ws.Range("K132").Value = 911772.48
$ws.Range("L132").Value = 16663.875
$ws.Range("M132").Value = -909242.48
$ws.Range("N132").Value = -21723.875
$ws.Range("H138").Value = 2893.7683
$ws.Range("I138").Value = 1133.7715
$ws.Range("J138").Value = 3920.4333
$ws.Range("K138").Value = 3401.3145
$ws.Range("L138").Value = 11761.2999
$ws.Range("M138").Value = 1738.6855
$ws.Range("N138").Value = -22041.2999
$ws.Range("H141").Value = 3240.8572
$ws.Range("I141").Value = 2337.7
$ws.Range("K141").Value = 7013.099999999999
$ws.Range("M141").Value = -1833.099999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3729.5273
$ws.Range("I2").Value = 1798.3489
$ws.Range("K2").Value = 1798.3489
$ws.Range("M2").Value = -1685.3489
$ws.Range("H32").Value = 2102522.2
$ws.Range("I32").Value = 5075.229
$ws.Range("J32").Value = 22238014
$ws.Range("K32").Value = 5075.229
$ws.Range("L32").Value = 22238014
$ws.Range("M32").Value = -4788.229
$ws.Range("N32").Value = -22238588
$ws.Range("H116").Value = 3729.5273
$ws.Range("I116").Value = 1798.3489
$ws.Range("K116").Value = 1798.3489
$ws.Range("M116").Value = 495.6511
$ws.Range("H132").Value = 1478301.1
$ws.Range("I132").Value = 2087760.2
$ws.Range("J132").Value = 123947.336
$ws.Range("K132").Value = 6263280.6
$ws.Range("L132").Value = 371842.008
$ws.Range("M132").Value = -6260750.6
$ws.Range("N132").Value = -376902.008

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3729.5273
$ws.Range("I3").Value = 1798.3489
$ws.Range("K3").Value = 1798.3489
$ws.Range("M3").Value = -1684.3489
$ws.Range("H99").Value = 7431.5396
$ws.Range("I99").Value = 7484.6816
$ws.Range("J99").Value = 7358.4688
$ws.Range("K99").Value = 7484.6816
$ws.Range("L99").Value = 7358.4688
$ws.Range("M99").Value = -5986.6816
$ws.Range("N99").Value = -10354.4688

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 13516813
$ws.Range("I16").Value = 16668166
$ws.Range("J16").Value = 11011.571
$ws.Range("K16").Value = 16668166
$ws.Range("L16").Value = 11011.571
$ws.Range("M16").Value = -16667879
$ws.Range("N16").Value = -11585.571
$ws.Range("H99").Value = 7941300
$ws.Range("I99").Value = 15877841
$ws.Range("J99").Value = 4758.4287
$ws.Range("K99").Value = 15877841
$ws.Range("L99").Value = 4758.4287
$ws.Range("M99").Value = -15876343
$ws.Range("N99").Value = -7754.4287
$ws.Range("H113").Value = 13516813
$ws.Range("I113").Value = 16668166
$ws.Range("J113").Value = 11011.571
$ws.Range("K113").Value = 16668166
$ws.Range("L113").Value = 11011.571
$ws.Range("M113").Value = -16665996
$ws.Range("N113").Value = -15351.571
$ws.Range("H126").Value = 7941300
$ws.Range("I126").Value = 15877841
$ws.Range("J126").Value = 4758.4287
$ws.Range("K126").Value = 47633523
$ws.Range("L126").Value = 14275.2861
$ws.Range("M126").Value = -47631053
$ws.Range("N126").Value = -19215.2861
$ws.Range("H134").Value = 6593.4
$ws.Range("I134").Value = 3096.9167
$ws.Range("J134").Value = 14222.091
$ws.Range("K134").Value = 9290.750100000001
$ws.Range("L134").Value = 42666.273
$ws.Range("M134").Value = -6755.750100000001
$ws.Range("N134").Value = -47736.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 581.625
$ws.Range("I8").Value = 581.625
$ws.Range("K8").Value = 1744.875
$ws.Range("M8").Value = -1605.875
$ws.Range("H37").Value = 112964.2
$ws.Range("J37").Value = 112964.2
$ws.Range("L37").Value = 338892.6
$ws.Range("N37").Value = -339116.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10009.8
$ws.Range("J80").Value = 11686.5
$ws.Range("L80").Value = 11686.5
$ws.Range("N80").Value = -13682.5
$ws.Range("H83").Value = 10009.8
$ws.Range("J83").Value = 11686.5
$ws.Range("L83").Value = 58432.5
$ws.Range("N83").Value = -68416.5
$ws.Range("H132").Value = 4576.294
$ws.Range("I132").Value = 4132.45
$ws.Range("J132").Value = 6190.273
$ws.Range("K132").Value = 12397.35
$ws.Range("L132").Value = 18570.819
$ws.Range("M132").Value = -9867.349999999999
$ws.Range("N132").Value = -23630.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3521.6365
$ws.Range("I100").Value = 5457
$ws.Range("K100").Value = 5457
$ws.Range("M100").Value = -4916
$ws.Range("H122").Value = 367577.6
$ws.Range("I122").Value = 448050.62
$ws.Range("J122").Value = 5448.9
$ws.Range("K122").Value = 1344151.86
$ws.Range("L122").Value = 16346.7
$ws.Range("M122").Value = -1341701.86
$ws.Range("N122").Value = -21246.7
$ws.Range("H132").Value = 4572.3906
$ws.Range("I132").Value = 3331.652
$ws.Range("K132").Value = 9994.956
$ws.Range("M132").Value = -7464.956

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10579.895
$ws.Range("I62").Value = 12387.5
$ws.Range("J62").Value = 10097.866
$ws.Range("K62").Value = 12387.5
$ws.Range("L62").Value = 10097.866
$ws.Range("M62").Value = -11763.5
$ws.Range("N62").Value = -11345.866
$ws.Range("H65").Value = 10579.895
$ws.Range("I65").Value = 12387.5
$ws.Range("J65").Value = 10097.866
$ws.Range("K65").Value = 61937.5
$ws.Range("L65").Value = 50489.33
$ws.Range("M65").Value = -58817.5
$ws.Range("N65").Value = -56729.33
$ws.Range("H136").Value = 8072770
$ws.Range("I136").Value = 13521452
$ws.Range("J136").Value = 8720.879999999999
$ws.Range("K136").Value = 40564356
$ws.Range("L136").Value = 26162.64
$ws.Range("M136").Value = -40561806
$ws.Range("N136").Value = -31262.64
